$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 330 is brand new; set every column (A through T).
$ws.Range("A330").Value = 11
$ws.Range("B330").Value = 'Vega Monumental Concepción'
$ws.Range("C330").Value = 'Bíobío'
$ws.Range("D330").Value = 44160
$ws.Range("E330").Value = 8
$ws.Range("F330").Value = 'Fruta'
$ws.Range("G330").Value = 100103
$ws.Range("H330").Value = 'Frutos de hueso (carozo)'
$ws.Range("I330").Value = 100103006
$ws.Range("J330").Value = 'Nectarín'
$ws.Range("K330").Value = 'Artic Star'
$ws.Range("L330").Value = 'Segunda'
$ws.Range("M330").Value = 100
$ws.Range("N330").Value = 16000
$ws.Range("O330").Value = 17000
$ws.Range("P330").Value = 16500
$ws.Range("Q330").Value = '$/caja 15 kilos empedrada'
$ws.Range("R330").Value = 'Región de O''Higgins'
$ws.Range("S330").Value = 1100
$ws.Range("T330").Value = 15
$ws.Range("D330").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Rows 277-329: shift values down one row (row r gets what used to be in row r-1).
# Row 277 receives the newly reported data point.
# Row 277
$ws.Range("D277").Value = 44637
$ws.Range("K277").Value = 'Artic Snow'
$ws.Range("L277").Value = 'Primera'
$ws.Range("M277").Value = 280
$ws.Range("N277").Value = 10000
$ws.Range("O277").Value = 11000
$ws.Range("P277").Value = 10464
$ws.Range("Q277").Value = '$/caja 16 kilos empedrada'
$ws.Range("S277").Value = 654
$ws.Range("T277").Value = 16

# Row 278
$ws.Range("D278").Value = 44208
$ws.Range("K278").Value = 'Super Queen'
$ws.Range("L278").Value = 'Especial'
$ws.Range("M278").Value = 100
$ws.Range("N278").Value = 14000
$ws.Range("O278").Value = 14000
$ws.Range("P278").Value = 14000
$ws.Range("Q278").Value = '$/caja 16 kilos empedrada'
$ws.Range("S278").Value = 875
$ws.Range("T278").Value = 16

# Row 279
$ws.Range("D279").Value = 44208
$ws.Range("K279").Value = 'Super Queen'
$ws.Range("L279").Value = 'Primera'
$ws.Range("M279").Value = 100
$ws.Range("N279").Value = 12000
$ws.Range("O279").Value = 12000
$ws.Range("P279").Value = 12000
$ws.Range("Q279").Value = '$/caja 16 kilos empedrada'
$ws.Range("S279").Value = 750
$ws.Range("T279").Value = 16

# Row 280
$ws.Range("D280").Value = 44208
$ws.Range("K280").Value = 'Super Queen'
$ws.Range("L280").Value = 'Segunda'
$ws.Range("M280").Value = 100
$ws.Range("N280").Value = 10000
$ws.Range("O280").Value = 10000
$ws.Range("P280").Value = 10000
$ws.Range("Q280").Value = '$/caja 16 kilos empedrada'
$ws.Range("S280").Value = 625
$ws.Range("T280").Value = 16

# Row 281
$ws.Range("D281").Value = 44530
$ws.Range("K281").Value = 'Big John'
$ws.Range("L281").Value = 'Primera'
$ws.Range("M281").Value = 200
$ws.Range("N281").Value = 15000
$ws.Range("O281").Value = 16000
$ws.Range("P281").Value = 15500
$ws.Range("Q281").Value = '$/caja 15 kilos empedrada'
$ws.Range("S281").Value = 1033
$ws.Range("T281").Value = 15

# Row 282
$ws.Range("D282").Value = 44530
$ws.Range("K282").Value = 'Big John'
$ws.Range("L282").Value = 'Segunda'
$ws.Range("M282").Value = 100
$ws.Range("N282").Value = 14000
$ws.Range("O282").Value = 14000
$ws.Range("P282").Value = 14000
$ws.Range("Q282").Value = '$/caja 15 kilos empedrada'
$ws.Range("S282").Value = 933
$ws.Range("T282").Value = 15

# Row 283
$ws.Range("D283").Value = 44617
$ws.Range("K283").Value = 'August Red'
$ws.Range("L283").Value = 'Especial'
$ws.Range("M283").Value = 50
$ws.Range("N283").Value = 12000
$ws.Range("O283").Value = 12000
$ws.Range("P283").Value = 12000
$ws.Range("Q283").Value = '$/bandeja 18 kilos granel'
$ws.Range("S283").Value = 667
$ws.Range("T283").Value = 18

# Row 284
$ws.Range("D284").Value = 44617
$ws.Range("K284").Value = 'August Red'
$ws.Range("L284").Value = 'Primera'
$ws.Range("M284").Value = 100
$ws.Range("N284").Value = 10000
$ws.Range("O284").Value = 10000
$ws.Range("P284").Value = 10000
$ws.Range("Q284").Value = '$/bandeja 18 kilos granel'
$ws.Range("S284").Value = 556
$ws.Range("T284").Value = 18

# Row 285
$ws.Range("D285").Value = 44617
$ws.Range("K285").Value = 'August Red'
$ws.Range("L285").Value = 'Segunda'
$ws.Range("M285").Value = 100
$ws.Range("N285").Value = 9000
$ws.Range("O285").Value = 9000
$ws.Range("P285").Value = 9000
$ws.Range("Q285").Value = '$/bandeja 18 kilos granel'
$ws.Range("S285").Value = 500
$ws.Range("T285").Value = 18

# Row 286
$ws.Range("D286").Value = 44264
$ws.Range("K286").Value = 'August Red'
$ws.Range("L286").Value = 'Primera'
$ws.Range("M286").Value = 200
$ws.Range("N286").Value = 11000
$ws.Range("O286").Value = 12000
$ws.Range("P286").Value = 11500
$ws.Range("Q286").Value = '$/caja 16 kilos empedrada'
$ws.Range("S286").Value = 719
$ws.Range("T286").Value = 16

# Row 287
$ws.Range("D287").Value = 44264
$ws.Range("K287").Value = 'August Red'
$ws.Range("L287").Value = 'Segunda'
$ws.Range("M287").Value = 100
$ws.Range("N287").Value = 10000
$ws.Range("O287").Value = 10000
$ws.Range("P287").Value = 10000
$ws.Range("Q287").Value = '$/caja 16 kilos empedrada'
$ws.Range("S287").Value = 625
$ws.Range("T287").Value = 16

# Row 288
$ws.Range("D288").Value = 44232
$ws.Range("K288").Value = 'Venus'
$ws.Range("L288").Value = 'Especial'
$ws.Range("M288").Value = 100
$ws.Range("N288").Value = 14000
$ws.Range("O288").Value = 14000
$ws.Range("P288").Value = 14000
$ws.Range("Q288").Value = '$/caja 16 kilos empedrada'
$ws.Range("S288").Value = 875
$ws.Range("T288").Value = 16

# Row 289
$ws.Range("D289").Value = 44232
$ws.Range("K289").Value = 'Venus'
$ws.Range("L289").Value = 'Primera'
$ws.Range("M289").Value = 100
$ws.Range("N289").Value = 12000
$ws.Range("O289").Value = 12000
$ws.Range("P289").Value = 12000
$ws.Range("Q289").Value = '$/caja 16 kilos empedrada'
$ws.Range("S289").Value = 750
$ws.Range("T289").Value = 16

# Row 290
$ws.Range("D290").Value = 44232
$ws.Range("K290").Value = 'Venus'
$ws.Range("L290").Value = 'Segunda'
$ws.Range("M290").Value = 100
$ws.Range("N290").Value = 10000
$ws.Range("O290").Value = 10000
$ws.Range("P290").Value = 10000
$ws.Range("Q290").Value = '$/caja 16 kilos empedrada'
$ws.Range("S290").Value = 625
$ws.Range("T290").Value = 16

# Row 291
$ws.Range("D291").Value = 44279
$ws.Range("K291").Value = 'Artic Mist'
$ws.Range("L291").Value = 'Primera'
$ws.Range("M291").Value = 100
$ws.Range("N291").Value = 12000
$ws.Range("O291").Value = 12000
$ws.Range("P291").Value = 12000
$ws.Range("Q291").Value = '$/caja 16 kilos empedrada'
$ws.Range("S291").Value = 750
$ws.Range("T291").Value = 16

# Row 292
$ws.Range("D292").Value = 44279
$ws.Range("K292").Value = 'Artic Mist'
$ws.Range("L292").Value = 'Segunda'
$ws.Range("M292").Value = 100
$ws.Range("N292").Value = 10000
$ws.Range("O292").Value = 10000
$ws.Range("P292").Value = 10000
$ws.Range("Q292").Value = '$/caja 16 kilos empedrada'
$ws.Range("S292").Value = 625
$ws.Range("T292").Value = 16

# Row 293
$ws.Range("D293").Value = 44572
$ws.Range("K293").Value = 'Red Diamond'
$ws.Range("L293").Value = 'Primera'
$ws.Range("M293").Value = 100
$ws.Range("N293").Value = 16000
$ws.Range("O293").Value = 16000
$ws.Range("P293").Value = 16000
$ws.Range("Q293").Value = '$/caja 15 kilos empedrada'
$ws.Range("S293").Value = 1067
$ws.Range("T293").Value = 15

# Row 294
$ws.Range("D294").Value = 44572
$ws.Range("K294").Value = 'Red Diamond'
$ws.Range("L294").Value = 'Segunda'
$ws.Range("M294").Value = 100
$ws.Range("N294").Value = 14000
$ws.Range("O294").Value = 14000
$ws.Range("P294").Value = 14000
$ws.Range("Q294").Value = '$/caja 15 kilos empedrada'
$ws.Range("S294").Value = 933
$ws.Range("T294").Value = 15

# Row 295
$ws.Range("D295").Value = 44257
$ws.Range("K295").Value = 'Artic Snow'
$ws.Range("L295").Value = 'Especial'
$ws.Range("M295").Value = 50
$ws.Range("N295").Value = 14000
$ws.Range("O295").Value = 14000
$ws.Range("P295").Value = 14000
$ws.Range("Q295").Value = '$/caja 16 kilos empedrada'
$ws.Range("S295").Value = 875
$ws.Range("T295").Value = 16

# Row 296
$ws.Range("D296").Value = 44257
$ws.Range("K296").Value = 'Artic Snow'
$ws.Range("L296").Value = 'Primera'
$ws.Range("M296").Value = 100
$ws.Range("N296").Value = 12000
$ws.Range("O296").Value = 12000
$ws.Range("P296").Value = 12000
$ws.Range("Q296").Value = '$/caja 16 kilos empedrada'
$ws.Range("S296").Value = 750
$ws.Range("T296").Value = 16

# Row 297
$ws.Range("D297").Value = 44257
$ws.Range("K297").Value = 'Artic Snow'
$ws.Range("L297").Value = 'Segunda'
$ws.Range("M297").Value = 100
$ws.Range("N297").Value = 10000
$ws.Range("O297").Value = 10000
$ws.Range("P297").Value = 10000
$ws.Range("Q297").Value = '$/caja 16 kilos empedrada'
$ws.Range("S297").Value = 625
$ws.Range("T297").Value = 16

# Row 298
$ws.Range("D298").Value = 44257
$ws.Range("K298").Value = 'August Red'
$ws.Range("L298").Value = 'Especial'
$ws.Range("M298").Value = 50
$ws.Range("N298").Value = 14000
$ws.Range("O298").Value = 14000
$ws.Range("P298").Value = 14000
$ws.Range("Q298").Value = '$/caja 16 kilos empedrada'
$ws.Range("S298").Value = 875
$ws.Range("T298").Value = 16

# Row 299
$ws.Range("D299").Value = 44257
$ws.Range("K299").Value = 'August Red'
$ws.Range("L299").Value = 'Primera'
$ws.Range("M299").Value = 100
$ws.Range("N299").Value = 12000
$ws.Range("O299").Value = 12000
$ws.Range("P299").Value = 12000
$ws.Range("Q299").Value = '$/caja 16 kilos empedrada'
$ws.Range("S299").Value = 750
$ws.Range("T299").Value = 16

# Row 300
$ws.Range("D300").Value = 44257
$ws.Range("K300").Value = 'August Red'
$ws.Range("L300").Value = 'Segunda'
$ws.Range("M300").Value = 100
$ws.Range("N300").Value = 10000
$ws.Range("O300").Value = 10000
$ws.Range("P300").Value = 10000
$ws.Range("Q300").Value = '$/caja 16 kilos empedrada'
$ws.Range("S300").Value = 625
$ws.Range("T300").Value = 16

# Row 301
$ws.Range("D301").Value = 44257
$ws.Range("K301").Value = 'Venus'
$ws.Range("L301").Value = 'Especial'
$ws.Range("M301").Value = 50
$ws.Range("N301").Value = 14000
$ws.Range("O301").Value = 14000
$ws.Range("P301").Value = 14000
$ws.Range("Q301").Value = '$/caja 16 kilos empedrada'
$ws.Range("S301").Value = 875
$ws.Range("T301").Value = 16

# Row 302
$ws.Range("D302").Value = 44257
$ws.Range("K302").Value = 'Venus'
$ws.Range("L302").Value = 'Primera'
$ws.Range("M302").Value = 100
$ws.Range("N302").Value = 12000
$ws.Range("O302").Value = 12000
$ws.Range("P302").Value = 12000
$ws.Range("Q302").Value = '$/caja 16 kilos empedrada'
$ws.Range("S302").Value = 750
$ws.Range("T302").Value = 16

# Row 303
$ws.Range("D303").Value = 44257
$ws.Range("K303").Value = 'Venus'
$ws.Range("L303").Value = 'Segunda'
$ws.Range("M303").Value = 100
$ws.Range("N303").Value = 10000
$ws.Range("O303").Value = 10000
$ws.Range("P303").Value = 10000
$ws.Range("Q303").Value = '$/caja 16 kilos empedrada'
$ws.Range("S303").Value = 625
$ws.Range("T303").Value = 16

# Row 304
$ws.Range("D304").Value = 44236
$ws.Range("K304").Value = 'Venus'
$ws.Range("L304").Value = 'Especial'
$ws.Range("M304").Value = 100
$ws.Range("N304").Value = 14000
$ws.Range("O304").Value = 14000
$ws.Range("P304").Value = 14000
$ws.Range("Q304").Value = '$/caja 16 kilos empedrada'
$ws.Range("S304").Value = 875
$ws.Range("T304").Value = 16

# Row 305
$ws.Range("D305").Value = 44236
$ws.Range("K305").Value = 'Venus'
$ws.Range("L305").Value = 'Primera'
$ws.Range("M305").Value = 100
$ws.Range("N305").Value = 12000
$ws.Range("O305").Value = 12000
$ws.Range("P305").Value = 12000
$ws.Range("Q305").Value = '$/caja 16 kilos empedrada'
$ws.Range("S305").Value = 750
$ws.Range("T305").Value = 16

# Row 306
$ws.Range("D306").Value = 44236
$ws.Range("K306").Value = 'Venus'
$ws.Range("L306").Value = 'Segunda'
$ws.Range("M306").Value = 100
$ws.Range("N306").Value = 10000
$ws.Range("O306").Value = 10000
$ws.Range("P306").Value = 10000
$ws.Range("Q306").Value = '$/caja 16 kilos empedrada'
$ws.Range("S306").Value = 625
$ws.Range("T306").Value = 16

# Row 307
$ws.Range("D307").Value = 44229
$ws.Range("K307").Value = 'Nectar Crest'
$ws.Range("L307").Value = 'Especial'
$ws.Range("M307").Value = 100
$ws.Range("N307").Value = 15000
$ws.Range("O307").Value = 15000
$ws.Range("P307").Value = 15000
$ws.Range("Q307").Value = '$/caja 16 kilos empedrada'
$ws.Range("S307").Value = 938
$ws.Range("T307").Value = 16

# Row 308
$ws.Range("D308").Value = 44229
$ws.Range("K308").Value = 'Nectar Crest'
$ws.Range("L308").Value = 'Primera'
$ws.Range("M308").Value = 100
$ws.Range("N308").Value = 13000
$ws.Range("O308").Value = 13000
$ws.Range("P308").Value = 13000
$ws.Range("Q308").Value = '$/caja 16 kilos empedrada'
$ws.Range("S308").Value = 812
$ws.Range("T308").Value = 16

# Row 309
$ws.Range("D309").Value = 44229
$ws.Range("K309").Value = 'Nectar Crest'
$ws.Range("L309").Value = 'Segunda'
$ws.Range("M309").Value = 100
$ws.Range("N309").Value = 11000
$ws.Range("O309").Value = 11000
$ws.Range("P309").Value = 11000
$ws.Range("Q309").Value = '$/caja 16 kilos empedrada'
$ws.Range("S309").Value = 688
$ws.Range("T309").Value = 16

# Row 310
$ws.Range("D310").Value = 44229
$ws.Range("K310").Value = 'Venus'
$ws.Range("L310").Value = 'Especial'
$ws.Range("M310").Value = 50
$ws.Range("N310").Value = 15000
$ws.Range("O310").Value = 15000
$ws.Range("P310").Value = 15000
$ws.Range("Q310").Value = '$/caja 16 kilos empedrada'
$ws.Range("S310").Value = 938
$ws.Range("T310").Value = 16

# Row 311
$ws.Range("D311").Value = 44229
$ws.Range("K311").Value = 'Venus'
$ws.Range("L311").Value = 'Primera'
$ws.Range("M311").Value = 100
$ws.Range("N311").Value = 13000
$ws.Range("O311").Value = 13000
$ws.Range("P311").Value = 13000
$ws.Range("Q311").Value = '$/caja 16 kilos empedrada'
$ws.Range("S311").Value = 812
$ws.Range("T311").Value = 16

# Row 312
$ws.Range("D312").Value = 44229
$ws.Range("K312").Value = 'Venus'
$ws.Range("L312").Value = 'Segunda'
$ws.Range("M312").Value = 100
$ws.Range("N312").Value = 11000
$ws.Range("O312").Value = 11000
$ws.Range("P312").Value = 11000
$ws.Range("Q312").Value = '$/caja 16 kilos empedrada'
$ws.Range("S312").Value = 688
$ws.Range("T312").Value = 16

# Row 313
$ws.Range("D313").Value = 44615
$ws.Range("K313").Value = 'August Red'
$ws.Range("L313").Value = 'Especial'
$ws.Range("M313").Value = 50
$ws.Range("N313").Value = 12000
$ws.Range("O313").Value = 12000
$ws.Range("P313").Value = 12000
$ws.Range("Q313").Value = '$/caja 16 kilos empedrada'
$ws.Range("S313").Value = 750
$ws.Range("T313").Value = 16

# Row 314
$ws.Range("D314").Value = 44615
$ws.Range("K314").Value = 'August Red'
$ws.Range("L314").Value = 'Primera'
$ws.Range("M314").Value = 100
$ws.Range("N314").Value = 10000
$ws.Range("O314").Value = 10000
$ws.Range("P314").Value = 10000
$ws.Range("Q314").Value = '$/caja 16 kilos empedrada'
$ws.Range("S314").Value = 625
$ws.Range("T314").Value = 16

# Row 315
$ws.Range("D315").Value = 44615
$ws.Range("K315").Value = 'August Red'
$ws.Range("L315").Value = 'Segunda'
$ws.Range("M315").Value = 100
$ws.Range("N315").Value = 8000
$ws.Range("O315").Value = 8000
$ws.Range("P315").Value = 8000
$ws.Range("Q315").Value = '$/caja 16 kilos empedrada'
$ws.Range("S315").Value = 500
$ws.Range("T315").Value = 16

# Row 316
$ws.Range("D316").Value = 44615
$ws.Range("K316").Value = 'June Pearl'
$ws.Range("L316").Value = 'Especial'
$ws.Range("M316").Value = 50
$ws.Range("N316").Value = 12000
$ws.Range("O316").Value = 12000
$ws.Range("P316").Value = 12000
$ws.Range("Q316").Value = '$/caja 16 kilos empedrada'
$ws.Range("S316").Value = 750
$ws.Range("T316").Value = 16

# Row 317
$ws.Range("D317").Value = 44615
$ws.Range("K317").Value = 'June Pearl'
$ws.Range("L317").Value = 'Primera'
$ws.Range("M317").Value = 100
$ws.Range("N317").Value = 10000
$ws.Range("O317").Value = 10000
$ws.Range("P317").Value = 10000
$ws.Range("Q317").Value = '$/caja 16 kilos empedrada'
$ws.Range("S317").Value = 625
$ws.Range("T317").Value = 16

# Row 318
$ws.Range("D318").Value = 44615
$ws.Range("K318").Value = 'June Pearl'
$ws.Range("L318").Value = 'Segunda'
$ws.Range("M318").Value = 100
$ws.Range("N318").Value = 8000
$ws.Range("O318").Value = 8000
$ws.Range("P318").Value = 8000
$ws.Range("Q318").Value = '$/caja 16 kilos empedrada'
$ws.Range("S318").Value = 500
$ws.Range("T318").Value = 16

# Row 319
$ws.Range("D319").Value = 44167
$ws.Range("K319").Value = 'Artic Star'
$ws.Range("L319").Value = 'Primera'
$ws.Range("M319").Value = 100
$ws.Range("N319").Value = 16000
$ws.Range("O319").Value = 17000
$ws.Range("P319").Value = 16500
$ws.Range("Q319").Value = '$/caja 16 kilos empedrada'
$ws.Range("S319").Value = 1031
$ws.Range("T319").Value = 16

# Row 320
$ws.Range("D320").Value = 44167
$ws.Range("K320").Value = 'Artic Star'
$ws.Range("L320").Value = 'Segunda'
$ws.Range("M320").Value = 50
$ws.Range("N320").Value = 15000
$ws.Range("O320").Value = 15000
$ws.Range("P320").Value = 15000
$ws.Range("Q320").Value = '$/caja 16 kilos empedrada'
$ws.Range("S320").Value = 938
$ws.Range("T320").Value = 16

# Row 321
$ws.Range("D321").Value = 44258
$ws.Range("K321").Value = 'Artic Snow'
$ws.Range("L321").Value = 'Primera'
$ws.Range("M321").Value = 100
$ws.Range("N321").Value = 12000
$ws.Range("O321").Value = 12000
$ws.Range("P321").Value = 12000
$ws.Range("Q321").Value = '$/caja 16 kilos empedrada'
$ws.Range("S321").Value = 750
$ws.Range("T321").Value = 16

# Row 322
$ws.Range("D322").Value = 44258
$ws.Range("K322").Value = 'Artic Snow'
$ws.Range("L322").Value = 'Segunda'
$ws.Range("M322").Value = 100
$ws.Range("N322").Value = 10000
$ws.Range("O322").Value = 10000
$ws.Range("P322").Value = 10000
$ws.Range("Q322").Value = '$/caja 16 kilos empedrada'
$ws.Range("S322").Value = 625
$ws.Range("T322").Value = 16

# Row 323
$ws.Range("D323").Value = 44285
$ws.Range("K323").Value = 'June Pearl'
$ws.Range("L323").Value = 'Primera'
$ws.Range("M323").Value = 200
$ws.Range("N323").Value = 11000
$ws.Range("O323").Value = 12000
$ws.Range("P323").Value = 11500
$ws.Range("Q323").Value = '$/caja 16 kilos empedrada'
$ws.Range("S323").Value = 719
$ws.Range("T323").Value = 16

# Row 324
$ws.Range("D324").Value = 44285
$ws.Range("K324").Value = 'June Pearl'
$ws.Range("L324").Value = 'Segunda'
$ws.Range("M324").Value = 100
$ws.Range("N324").Value = 10000
$ws.Range("O324").Value = 10000
$ws.Range("P324").Value = 10000
$ws.Range("Q324").Value = '$/caja 16 kilos empedrada'
$ws.Range("S324").Value = 625
$ws.Range("T324").Value = 16

# Row 325
$ws.Range("D325").Value = 44595
$ws.Range("K325").Value = 'Red Diamond'
$ws.Range("L325").Value = 'Primera'
$ws.Range("M325").Value = 220
$ws.Range("N325").Value = 9500
$ws.Range("O325").Value = 10000
$ws.Range("P325").Value = 9773
$ws.Range("Q325").Value = '$/caja 16 kilos empedrada'
$ws.Range("S325").Value = 611
$ws.Range("T325").Value = 16

# Row 326
$ws.Range("D326").Value = 44628
$ws.Range("K326").Value = 'June Pearl'
$ws.Range("L326").Value = 'Especial'
$ws.Range("M326").Value = 50
$ws.Range("N326").Value = 15000
$ws.Range("O326").Value = 15000
$ws.Range("P326").Value = 15000
$ws.Range("Q326").Value = '$/bandeja 18 kilos granel'
$ws.Range("S326").Value = 833
$ws.Range("T326").Value = 18

# Row 327
$ws.Range("D327").Value = 44628
$ws.Range("K327").Value = 'June Pearl'
$ws.Range("L327").Value = 'Primera'
$ws.Range("M327").Value = 100
$ws.Range("N327").Value = 13000
$ws.Range("O327").Value = 13000
$ws.Range("P327").Value = 13000
$ws.Range("Q327").Value = '$/bandeja 18 kilos granel'
$ws.Range("S327").Value = 722
$ws.Range("T327").Value = 18

# Row 328
$ws.Range("D328").Value = 44628
$ws.Range("K328").Value = 'June Pearl'
$ws.Range("L328").Value = 'Segunda'
$ws.Range("M328").Value = 100
$ws.Range("N328").Value = 11000
$ws.Range("O328").Value = 11000
$ws.Range("P328").Value = 11000
$ws.Range("Q328").Value = '$/bandeja 18 kilos granel'
$ws.Range("S328").Value = 611
$ws.Range("T328").Value = 18

# Row 329
$ws.Range("D329").Value = 44160
$ws.Range("K329").Value = 'Artic Star'
$ws.Range("L329").Value = 'Primera'
$ws.Range("M329").Value = 100
$ws.Range("N329").Value = 19000
$ws.Range("O329").Value = 20000
$ws.Range("P329").Value = 19500
$ws.Range("Q329").Value = '$/caja 15 kilos empedrada'
$ws.Range("S329").Value = 1300
$ws.Range("T329").Value = 15
